$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Support")
$ws.Activate()

# Rename header for column B ("Program Funding" -> "Program Funding Sources")
$ws.Range("B1").Value = "Program Funding Sources"

# Update the funding-source descriptions (write in an order that reproduces
# the original shared-string insertion order: City, GGEE, District, Local,
# then the unchanged / repeated values).
$ws.Range("B4").Value = "City Funded, State AI, GGEE"
$ws.Range("B5").Value = "GGEE Donor Funding, , State AI"
$ws.Range("B6").Value = "District Funding, State AI, GGEE"
$ws.Range("B2").Value = "Local Funding, State AI, GGEE"

$ws.Range("B3").Value = "UF Donor Funding"
$ws.Range("B7").Value = "District Funding, State AI, GGEE"
$ws.Range("B8").Value = "District Funding, State AI, GGEE"
$ws.Range("B9").Value = "District Funding, State AI, GGEE"
$ws.Range("B10").Value = "District Funding, State AI, GGEE"

# Widen column B to fit the new, longer text
$ws.Columns.Item(2).ColumnWidth = 31.67

# Clear the (no-op) fill formatting left over on the header/number column so
# the cells fall back to the shared default style
for ($i = 2; $i -le 10; $i++) {
    $ws.Cells.Item($i, 3).Interior.Pattern = -4142
}
$ws.Range("C1").Interior.Pattern = -4142

# Leave the selection on B2, matching the saved view state
$ws.Range("B2").Select()
